$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format column D as Text so numeric-looking strings (e.g. "194.93") are
# stored as literal text instead of being auto-converted to numbers by Excel,
# matching the source workbook where these are inline strings.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "69.314.17"
$ws.Range("E2").Value = "  -1.20%  "

# Row 3
$ws.Range("D3").Value = "3.532.35"
$ws.Range("E3").Value = "  -1.98%  "

# Row 4
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").Value = "194.93"
$ws.Range("E5").Value = "  -0.91%  "

# Row 6
$ws.Range("D6").Value = "581.81"
$ws.Range("E6").Value = "  -3.76%  "

# Row 7
$ws.Range("D7").Value = "0.608"
$ws.Range("E7").Value = "  -2.77%  "

# Row 8
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("E9").Value = "  -1.52%  "

# Row 10
$ws.Range("D10").Value = "0.628"
$ws.Range("E10").Value = "  -3.13%  "

# Row 11
$ws.Range("D11").Value = "51.53"
$ws.Range("E11").Value = "  -4.26%  "

# Row 12
$ws.Range("D12").Value = "0.0000285"
$ws.Range("E12").Value = "  -6.23%  "

# Row 13
$ws.Range("D13").Value = "9.19"
$ws.Range("E13").Value = "  -3.97%  "

# Row 14
$ws.Range("D14").Value = "4.093.86"
$ws.Range("E14").Value = "  -1.99%  "

# Row 15
$ws.Range("D15").Value = "664.78"
$ws.Range("E15").Value = "  +11.93%  "

# Row 16
$ws.Range("D16").Value = "69.405.15"
$ws.Range("E16").Value = "  -1.28%  "

# Row 17
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.538.02"
$ws.Range("E17").Value = "  -1.88%  "

# Row 18
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "12.48"
$ws.Range("E18").Value = "  -4.57%  "

# Row 19
$ws.Range("E19").Value = "  -0.76%  "

# Row 20
$ws.Range("D20").Value = "18.41"
$ws.Range("E20").Value = "  -3.70%  "

# Row 21
$ws.Range("D21").Value = "0.962"
$ws.Range("E21").Value = "  -3.34%  "

# Row 22
$ws.Range("D22").Value = "18.22"
$ws.Range("E22").Value = "  +2.88%  "

# Row 23
$ws.Range("D23").Value = "5.31"
$ws.Range("E23").Value = "  +2.86%  "

# Row 24
$ws.Range("D24").Value = "104.24"
$ws.Range("E24").Value = "  +2.66%  "

# Row 25
$ws.Range("D25").Value = "4.36"
$ws.Range("E25").Value = "  -5.43%  "

# Row 26
$ws.Range("E26").Value = "  -4.16%  "

# Row 27
$ws.Range("D27").Value = "10.12"
$ws.Range("E27").Value = "  -5.77%  "

# Row 28
$ws.Range("D28").Value = "9.56"
$ws.Range("E28").Value = "  -0.26%  "

# Row 29
$ws.Range("D29").Value = "32.97"
$ws.Range("E29").Value = "  -2.48%  "

# Row 30
$ws.Range("D30").Value = "4.36"
$ws.Range("E30").Value = "  -7.48%  "

# Row 31
$ws.Range("D31").Value = "6.74"
$ws.Range("E31").Value = "  -5.48%  "

# Row 32
$ws.Range("D32").Value = "11.71"
$ws.Range("E32").Value = "  -4.75%  "

# Row 33
$ws.Range("D33").Value = "0.110"
$ws.Range("E33").Value = "  -5.50%  "

# Row 34
$ws.Range("D34").Value = "61.86"
$ws.Range("E34").Value = "  -2.22%  "

# Row 35
$ws.Range("D35").Value = "3.779.76"
$ws.Range("E35").Value = "  -4.05%  "

# Row 36
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.02%  "

# Row 37
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0810"
$ws.Range("E37").Value = "  -8.55%  "

# Row 38
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "3.68"
$ws.Range("E38").Value = "  +4.28%  "

# Row 39
$ws.Range("D39").Value = "500.38"
$ws.Range("E39").Value = "  -4.71%  "

# Row 40
$ws.Range("D40").Value = "2.90"
$ws.Range("E40").Value = "  -6.54%  "

# Row 41
$ws.Range("D41").Value = "0.369"
$ws.Range("E41").Value = "  -5.33%  "

# Row 42
$ws.Range("E42").Value = "  +0.34%  "

# Row 43
$ws.Range("D43").Value = "34.45"
$ws.Range("E43").Value = "  -6.48%  "

# Row 44
$ws.Range("E44").Value = "  -1.65%  "

# Row 45
$ws.Range("D45").Value = "3.38"
$ws.Range("E45").Value = "  -1.59%  "

# Row 46
$ws.Range("D46").Value = "2.84"

# Row 47
$ws.Range("E47").Value = "  -3.12%  "

# Row 48
$ws.Range("E48").Value = "  -0.17%  "

# Row 49
$ws.Range("D49").Value = "8.28"
$ws.Range("E49").Value = "  -3.92%  "

# Row 50
$ws.Range("D50").Value = "1.75"
$ws.Range("E50").Value = "  +18.23%  "

# Row 51
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D51").Value = "2.69"
$ws.Range("E51").Value = "  +63.28%  "

# Restore the default (Normal) style on column D so no lingering text-format
# style is left on these cells (matches original workbook styling).
$ws.Range("D2:D51").Style = "Normal"
